# Commit: "Add shortcut to rename files/directories"
#
# Adds a new row to the "Terminal Commands" sheet documenting the `mv`
# rename shortcut, and clears the (visually-identical, font-only) style
# that the whole column A used to carry so the new row's plain cells
# match the rest of the column.

$wb = $excel.ActiveWorkbook

$gitSheet = $wb.Worksheets.Item("Git Bash Commanda")
$termSheet = $wb.Worksheets.Item("Terminal Commands")

# --- New content: "mv old_file_name new_file_name" / "rename file/directory" ---
$termSheet.Range("A19").Value = "mv old_file_name new_file_name"
$termSheet.Range("B19").Value = "rename file/directory"

# Move the active selection to the newly added row, like the author did.
[void]$termSheet.Range("B19").Select()

# --- Strip the redundant "applyFont" style (cellXfs idx 5) ---
# It only re-applies the default font, so clearing it back to "Normal"
# is visually a no-op but matches the cleaned-up style table.
$gitSheet.Range("A7").Style = "Normal"
$gitSheet.Range("A8").Style = "Normal"

for ($r = 2; $r -le 19; $r++) {
    $termSheet.Range("A$r").Style = "Normal"
}

# --- Column A on "Terminal Commands" widens to fit the longer command text ---
$termSheet.Columns.Item(1).ColumnWidth = 30.1666666666667
